# Apply cryptos list update (Tue Oct  8 17:44:41 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.216.63'
$ws.Range('E2').Value = '  -2.30%  '
$ws.Range('D3').Value = '2.437.55'
$ws.Range('E3').Value = '  -1.50%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = "'579.31"
$ws.Range('E5').Value = '  +0.55%  '
$ws.Range('D6').Value = "'142.71"
$ws.Range('E6').Value = '  -4.22%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  -2.52%  '
$ws.Range('D9').Value = '2.436.22'
$ws.Range('E9').Value = '  -1.43%  '
$ws.Range('E10').Value = '  -4.99%  '
$ws.Range('E11').Value = '  +1.35%  '
$ws.Range('E12').Value = '  -2.10%  '
$ws.Range('E13').Value = '  -3.80%  '
$ws.Range('D14').Value = "'26.33"
$ws.Range('E14').Value = '  -3.46%  '
$ws.Range('E15').Value = '  -5.57%  '
$ws.Range('D16').Value = '2.868.67'
$ws.Range('E16').Value = '  -1.98%  '
$ws.Range('D17').Value = '62.265.75'
$ws.Range('E17').Value = '  -1.89%  '
$ws.Range('D18').Value = '2.428.25'
$ws.Range('E18').Value = '  -2.13%  '
$ws.Range('D19').Value = "'10.92"
$ws.Range('E19').Value = '  -4.85%  '
$ws.Range('E20').Value = '  -4.84%  '
$ws.Range('D21').Value = "'328.87"
$ws.Range('E21').Value = '  -0.66%  '
$ws.Range('E22').Value = '  -3.03%  '
$ws.Range('E23').Value = '  -8.30%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').Value = "'65.54"
$ws.Range('E25').Value = '  -0.74%  '
$ws.Range('B26').Value = 'Bittensor'
$ws.Range('C26').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D26').Value = "'629.68"
$ws.Range('E26').Value = '  -0.23%  '
$ws.Range('B27').Value = 'Aptos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D27').Value = "'9.22"
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('E29').Value = '  -10.08%  '
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('E31').Value = '  -7.50%  '
$ws.Range('D32').Value = "'8.00"
$ws.Range('E32').Value = '  -5.07%  '
$ws.Range('D33').Value = "'0.142"
$ws.Range('E33').Value = '  -1.78%  '
$ws.Range('E34').Value = '  -1.94%  '
$ws.Range('D35').Value = "'4.93"
$ws.Range('E35').Value = '  -6.57%  '
$ws.Range('E36').Value = '  +0.23%  '
$ws.Range('E37').Value = '  -7.99%  '
$ws.Range('E38').Value = '  -2.46%  '
$ws.Range('D39').Value = "'149.30"
$ws.Range('E39').Value = '  +1.32%  '
$ws.Range('D40').Value = "'18.29"
$ws.Range('E40').Value = '  -3.29%  '
$ws.Range('D41').Value = "'5.21"
$ws.Range('E41').Value = '  -5.56%  '
$ws.Range('D42').Value = "'1.76"
$ws.Range('E42').Value = '  -3.52%  '
$ws.Range('D43').Value = "'42.71"
$ws.Range('E43').Value = '  +1.60%  '
$ws.Range('D45').Value = "'2.46"
$ws.Range('E45').Value = '  -10.68%  '
$ws.Range('D46').Value = "'142.59"
$ws.Range('E46').Value = '  -5.64%  '
$ws.Range('D47').Value = "'3.63"
$ws.Range('E47').Value = '  -4.10%  '
$ws.Range('E48').Value = '  -3.96%  '
$ws.Range('E49').Value = '  -1.91%  '
$ws.Range('E50').Value = '  -9.47%  '
$ws.Range('E51').Value = '  +2.48%  '
